$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$genStyle = $ws.Range("B2").Style

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "67.214.57"
$ws.Range("D2").Style = $genStyle
$ws.Range("E2").Value = "  -1.68%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.628.13"
$ws.Range("D3").Style = $genStyle
$ws.Range("E3").Value = "  -1.92%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("D4").Style = $genStyle
$ws.Range("E4").Value = "  +0.07%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "586.73"
$ws.Range("D5").Style = $genStyle
$ws.Range("E5").Value = "  -2.41%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "187.15"
$ws.Range("D6").Style = $genStyle
$ws.Range("E6").Value = "  +2.97%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.613"
$ws.Range("D7").Style = $genStyle
$ws.Range("E7").Value = "  -2.87%  "
$ws.Range("E8").Value = "  +0.31%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.683"
$ws.Range("D9").Style = $genStyle
$ws.Range("E9").Value = "  -4.75%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.149"
$ws.Range("D10").Style = $genStyle
$ws.Range("E10").Value = "  -8.62%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "55.09"
$ws.Range("D11").Style = $genStyle
$ws.Range("E11").Value = "  -2.73%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0000259"
$ws.Range("D12").Style = $genStyle
$ws.Range("E12").Value = "  -11.07%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "10.03"
$ws.Range("D13").Style = $genStyle
$ws.Range("E13").Value = "  -3.75%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "4.223.51"
$ws.Range("D14").Style = $genStyle
$ws.Range("E14").Value = "  -1.20%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.640.02"
$ws.Range("D15").Style = $genStyle
$ws.Range("E15").Value = "  -1.39%  "
$ws.Range("E16").Value = "  -0.10%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "18.53"
$ws.Range("D17").Style = $genStyle
$ws.Range("E17").Value = "  -4.36%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "67.106.97"
$ws.Range("D18").Style = $genStyle
$ws.Range("E18").Value = "  -1.44%  "
$ws.Range("E19").Value = "  -3.41%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "12.34"
$ws.Range("D20").Style = $genStyle
$ws.Range("E20").Value = "  -3.98%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "399.55"
$ws.Range("D21").Style = $genStyle
$ws.Range("E21").Value = "  -2.40%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.35"
$ws.Range("D22").Style = $genStyle
$ws.Range("E22").Value = "  -5.78%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "86.48"
$ws.Range("D23").Style = $genStyle
$ws.Range("E23").Value = "  -2.96%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.88"
$ws.Range("D24").Style = $genStyle
$ws.Range("E24").Value = "  -4.85%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "12.28"
$ws.Range("D25").Style = $genStyle
$ws.Range("E25").Value = "  -4.37%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "10.59"
$ws.Range("D26").Style = $genStyle
$ws.Range("E26").Value = "  -2.89%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "6.05"
$ws.Range("D27").Style = $genStyle
$ws.Range("E27").Value = "  +0.04%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "3.61"
$ws.Range("D28").Style = $genStyle
$ws.Range("E28").Value = "  -6.99%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "9.12"
$ws.Range("D29").Style = $genStyle
$ws.Range("E29").Value = "  -3.76%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "31.45"
$ws.Range("D30").Style = $genStyle
$ws.Range("E30").Value = "  -4.40%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "6.92"
$ws.Range("D31").Style = $genStyle
$ws.Range("E31").Value = "  -5.16%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "66.39"
$ws.Range("D32").Style = $genStyle
$ws.Range("E32").Value = "  +2.95%  "
$ws.Range("E33").Value = "  -3.81%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "43.13"
$ws.Range("D34").Style = $genStyle
$ws.Range("E34").Value = "  -0.88%  "
$ws.Range("E35").Value = "  -3.44%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "583.27"
$ws.Range("D36").Style = $genStyle
$ws.Range("E36").Value = "  -2.34%  "
$ws.Range("E37").Value = "  +0.01%  "
$ws.Range("E38").Value = "  +0.20%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.381"
$ws.Range("D39").Style = $genStyle
$ws.Range("E39").Value = "  -5.00%  "
$ws.Range("B40").Value = "PEPE"
$ws.Range("C40").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0₃0740"
$ws.Range("D40").Style = $genStyle
$ws.Range("E40").Value = "  -16.91%  "
$ws.Range("B41").Value = "Kaspa"
$ws.Range("C41").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.134"
$ws.Range("D41").Style = $genStyle
$ws.Range("E41").Value = "  -1.36%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.82"
$ws.Range("D42").Style = $genStyle
$ws.Range("E42").Value = "  -6.30%  "
$ws.Range("E43").Value = "  -4.98%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "3.16"
$ws.Range("D44").Style = $genStyle
$ws.Range("E44").Value = "  +0.07%  "
$ws.Range("B45").Value = "Stellar"
$ws.Range("C45").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.133"
$ws.Range("D45").Style = $genStyle
$ws.Range("E45").Value = "  -1.12%  "
$ws.Range("B46").Value = "Fetch.AI"
$ws.Range("C46").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.42"
$ws.Range("D46").Style = $genStyle
$ws.Range("E46").Value = "  -13.07%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.709.88"
$ws.Range("D47").Style = $genStyle
$ws.Range("E47").Value = "  -1.78%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "140.97"
$ws.Range("D48").Style = $genStyle
$ws.Range("E48").Value = "  -0.70%  "
$ws.Range("E49").Value = "  -7.19%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.55"
$ws.Range("D50").Style = $genStyle
$ws.Range("E50").Value = "  -6.51%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.62"
$ws.Range("D51").Style = $genStyle
$ws.Range("E51").Value = "  -5.56%  "
